$d = $word.ActiveDocument

# The edit (per sam2.txt / sam2.docx diff):
#   - para 1 keeps "AJLSJDLASLDJALSJD" but gets a second run "rrdrhrdh"
#     appended, with the whole paragraph wrapped in spellcheck
#     proofErr markers (spellStart/spellEnd) - this is what Word does
#     once it has finished "checking" a misspelled word.
#   - a new blank paragraph
#   - a new paragraph "Uthtddhd" (spell-checked, so proofErr-wrapped)
#   - a new paragraph "Vhjgvjgvjgjg" (spell-checked, so proofErr-wrapped)
#   - a new paragraph "khkhkhlkhkhkhk" - this is the word currently being
#     typed, so it has not been through the proofErr pass yet.
#
# InsertXML replaces the exact contents of the range it's called on, so
# we collapse to the end of paragraph 1 and hand it the complete OOXML
# for the final five paragraphs in one shot (this also avoids leaving
# behind stray empty runs that Paragraph/InsertParagraphAfter-style
# incremental edits would create for the blank paragraph).

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newContent =
  "<w:p $wNs>" +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>AJLSJDLASLDJALSJD</w:t></w:r>' +
    '<w:r><w:t>rrdrhrdh</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>' +
  "<w:p $wNs/>" +
  "<w:p $wNs>" +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Uthtddhd</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>' +
  "<w:p $wNs>" +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Vhjgvjgvjgjg</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>' +
  "<w:p $wNs>" +
    '<w:r><w:t>khkhkhlkhkhkhk</w:t></w:r>' +
  '</w:p>'

$target = $d.Paragraphs(1).Range
$target.Collapse(0)
[void]$target.InsertXML($newContent)
